$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark that sits after
#    "Chi phi dao tao 3000000".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Delete the two trailing bulleted list items under "Quan ly ma nguon":
#    "So do cac branch duoc tao ra" and "So dong lenh cua du an".
$rng = $d.Content
$found = $rng.Find.Execute("Sơ đồ các branch được tạo ra")
if ($found) {
    $para = $rng.Paragraphs(1)
    $para.Range.Delete()
}

$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Số dòng lệnh của dự án")
if ($found2) {
    $para2 = $rng2.Paragraphs(1)
    $para2.Range.Delete()
}

# 3. Re-add a "_GoBack" bookmark immediately after "thay doi" (the last
#    edit point), collapsed right before the paragraph mark. A directly
#    collapsed Range placed exactly at a paragraph's text-end offset is
#    mis-anchored by this runtime, so insert a throwaway marker
#    character, anchor the bookmark to it, then delete the marker -
#    the bookmark collapses back to the correct position.
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("dòng lệnh bị thay đổi")
if ($found3) {
    $endPos = $rng3.End
    $rng3.Collapse(0)
    $rng3.InsertAfter("X")
    $d.Bookmarks.Add("_GoBack", $rng3)
    $delRange = $d.Range($endPos, $endPos + 1)
    $delRange.Delete()
}
